{"js": "// Update the \"two-digit \u00f7 one-digit\" practice table cells with new\n// dividend/divisor/quotient/remainder values.\n// Each entry is [rowIndex, colIndex, oldText, newText] addressing a\n// specific cell of the (single) table in the document body.\nconst replacements = [\n  [0, 0, \"33\u00f79=3, 6\", \"72\u00f75=14, 2\"],\n  [0, 1, \"67\u00f78=8, 3\", \"12\u00f77=1, 5\"],\n  [0, 2, \"10\u00f78=1, 2\", \"17\u00f79=1, 8\"],\n  [0, 3, \"20\u00f75=4, 0\", \"63\u00f72=31, 1\"],\n  [0, 4, \"23\u00f77=3, 2\", \"57\u00f72=28, 1\"],\n  [4, 0, \"35\u00f72=17, 1\", \"62\u00f77=8, 6\"],\n  [4, 1, \"77\u00f79=8, 5\", \"76\u00f76=12, 4\"],\n  [4, 2, \"83\u00f75=16, 3\", \"95\u00f74=23, 3\"],\n  [4, 3, \"19\u00f77=2, 5\", \"33\u00f73=11, 0\"],\n  [4, 4, \"25\u00f74=6, 1\", \"67\u00f74=16, 3\"],\n  [8, 0, \"12\u00f72=6, 0\", \"96\u00f76=16, 0\"],\n  [8, 1, \"36\u00f72=18, 0\", \"38\u00f72=19, 0\"],\n  [8, 2, \"84\u00f76=14, 0\", \"78\u00f72=39, 0\"],\n  [8, 3, \"16\u00f77=2, 2\", \"45\u00f77=6, 3\"],\n  [8, 4, \"37\u00f76=6, 1\", \"97\u00f73=32, 1\"],\n  [12, 0, \"67\u00f75=13, 2\", \"28\u00f75=5, 3\"],\n  [12, 1, \"46\u00f73=15, 1\", \"72\u00f78=9, 0\"],\n  [12, 2, \"89\u00f76=14, 5\", \"13\u00f75=2, 3\"],\n  [12, 3, \"51\u00f78=6, 3\", \"44\u00f76=7, 2\"],\n  [12, 4, \"90\u00f75=18, 0\", \"65\u00f73=21, 2\"],\n  [16, 0, \"68\u00f75=13, 3\", \"60\u00f72=30, 0\"],\n  [16, 1, \"72\u00f78=9, 0\", \"60\u00f77=8, 4\"],\n  [16, 2, \"58\u00f79=6, 4\", \"52\u00f79=5, 7\"],\n  [16, 3, \"23\u00f76=3, 5\", \"47\u00f75=9, 2\"],\n  [16, 4, \"31\u00f74=7, 3\", \"42\u00f76=7, 0\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body.\");\n}\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${oldText}\" in cell (${row}, ${col}), found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the \"two-digit / one-digit\" division practice table: each\n# entry below is (row, column, oldExpression, newExpression) using\n# 1-based Word COM Table.Cell(row, column) addressing.\n#\n# We intentionally avoid Find/Replace here: this host's Range.Find is\n# not bounded to the range it was created from (it can match text\n# anywhere in the document), which is unsafe for this table because\n# several of the new expressions equal other cells' old expressions\n# elsewhere in the grid (e.g. \"72\u00f78=9, 0\" is both a new value for one\n# cell and the old value of a different cell). Instead we address each\n# cell directly and overwrite its Range.Text, trimming the trailing\n# end-of-cell marker first so the cell's paragraph/run formatting is\n# preserved.\n$replacements = @(\n  @(1, 1, '33\u00f79=3, 6', '72\u00f75=14, 2'),\n  @(1, 2, '67\u00f78=8, 3', '12\u00f77=1, 5'),\n  @(1, 3, '10\u00f78=1, 2', '17\u00f79=1, 8'),\n  @(1, 4, '20\u00f75=4, 0', '63\u00f72=31, 1'),\n  @(1, 5, '23\u00f77=3, 2', '57\u00f72=28, 1'),\n  @(5, 1, '35\u00f72=17, 1', '62\u00f77=8, 6'),\n  @(5, 2, '77\u00f79=8, 5', '76\u00f76=12, 4'),\n  @(5, 3, '83\u00f75=16, 3', '95\u00f74=23, 3'),\n  @(5, 4, '19\u00f77=2, 5', '33\u00f73=11, 0'),\n  @(5, 5, '25\u00f74=6, 1', '67\u00f74=16, 3'),\n  @(9, 1, '12\u00f72=6, 0', '96\u00f76=16, 0'),\n  @(9, 2, '36\u00f72=18, 0', '38\u00f72=19, 0'),\n  @(9, 3, '84\u00f76=14, 0', '78\u00f72=39, 0'),\n  @(9, 4, '16\u00f77=2, 2', '45\u00f77=6, 3'),\n  @(9, 5, '37\u00f76=6, 1', '97\u00f73=32, 1'),\n  @(13, 1, '67\u00f75=13, 2', '28\u00f75=5, 3'),\n  @(13, 2, '46\u00f73=15, 1', '72\u00f78=9, 0'),\n  @(13, 3, '89\u00f76=14, 5', '13\u00f75=2, 3'),\n  @(13, 4, '51\u00f78=6, 3', '44\u00f76=7, 2'),\n  @(13, 5, '90\u00f75=18, 0', '65\u00f73=21, 2'),\n  @(17, 1, '68\u00f75=13, 3', '60\u00f72=30, 0'),\n  @(17, 2, '72\u00f78=9, 0', '60\u00f77=8, 4'),\n  @(17, 3, '58\u00f79=6, 4', '52\u00f79=5, 7'),\n  @(17, 4, '23\u00f76=3, 5', '47\u00f75=9, 2'),\n  @(17, 5, '31\u00f74=7, 3', '42\u00f76=7, 0')\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nforeach ($entry in $replacements) {\n  $row = $entry[0]\n  $col = $entry[1]\n  $oldText = $entry[2]\n  $newText = $entry[3]\n\n  $cell = $t.Cell($row, $col)\n  $rng = $cell.Range\n  # Drop the trailing end-of-cell mark (and paragraph mark) from the\n  # range so only the visible text is considered/replaced.\n  $rng.MoveEnd(1, -1) | Out-Null\n\n  if ($rng.Text -ne $oldText) {\n    throw (\"Unexpected text in cell (\" + $row + \",\" + $col + \"): expected '\" + $oldText + \"' but found '\" + $rng.Text + \"'\")\n  }\n\n  $rng.Text = $newText\n}\n"}
